# Fixed a bug in FixSymbols
# Reorders the data rows (A2:F23) so that each row's symbol-id (column A)
# and its associated reel values end up on the corrected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(301, 6, 45, 30, 60, 45),
    @(601, 9, 60, 67, 60, 42),
    @(901, 16, 15, 45, 60, 60),
    @(902, 1, 0, 0, 0, 0),
    @(1001, 18, 30, 75, 60, 72),
    @(401, 9, 48, 67, 75, 45),
    @(201, 9, 30, 15, 45, 30),
    @(1202, 2, 10, 10, 10, 10),
    @(101, 9, 30, 15, 60, 15),
    @(1201, 2, 10, 10, 10, 10),
    @(801, 3, 67, 65, 52, 45),
    @(701, 3, 90, 45, 97, 15),
    @(1203, 3, 15, 15, 15, 15),
    @(501, 9, 52, 30, 75, 45),
    @(1101, 0, 15, 30, 30, 0),
    @(3, 0, 3, 3, 3, 3),
    @(502, 0, 4, 0, 0, 0),
    @(802, 0, 4, 5, 4, 0),
    @(2, 0, 2, 2, 2, 2),
    @(1, 0, 2, 2, 2, 2),
    @(602, 0, 0, 4, 0, 9),
    @(402, 0, 0, 4, 0, 0)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $values[$c]
    }
}
